$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.580.30'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '1.924.16'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('D4').Value = "'1.015"
$ws.Range('E4').Value = '  +0.72%  '
$ws.Range('D5').Value = "'326.74"
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').Value = "'1.013"
$ws.Range('E6').Value = '  +0.77%  '
$ws.Range('D7').Value = "'0.4818"
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('D8').Value = "'0.4052"
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('D9').Value = "'0.08203"
$ws.Range('E9').Value = '  +0.50%  '
$ws.Range('D10').Value = "'1.007"
$ws.Range('E10').Value = '  -0.52%  '
$ws.Range('D11').Value = "'23.69"
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.923.16'
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'6.094"
$ws.Range('E13').Value = '  +0.98%  '
$ws.Range('D14').Value = "'7.289"
$ws.Range('E14').Value = '  +1.44%  '
$ws.Range('D15').Value = "'91.57"
$ws.Range('E15').Value = '  +1.24%  '
$ws.Range('D16').Value = "'0.06879"
$ws.Range('E16').Value = '  +1.72%  '
$ws.Range('D17').Value = "'1.015"
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').Value = "'0.00001039"
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').Value = "'17.62"
$ws.Range('E19').Value = '  -0.27%  '
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('D21').Value = '29.577.24'
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').Value = "'5.669"
$ws.Range('E22').Value = '  +0.85%  '
$ws.Range('D23').Value = "'11.96"
$ws.Range('E23').Value = '  +1.64%  '
$ws.Range('D24').Value = "'2.185"
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').Value = '2.142.99'
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('D26').Value = "'156.24"
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('D27').Value = "'6.384"
$ws.Range('E27').Value = '  -1.83%  '
$ws.Range('D28').Value = "'20.02"
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('D29').Value = "'2.089"
$ws.Range('E29').Value = '  -1.25%  '
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('D31').Value = "'1.009"
$ws.Range('E31').Value = '  -1.86%  '
$ws.Range('D32').Value = "'0.09598"
$ws.Range('E32').Value = '  +0.74%  '
$ws.Range('D33').Value = "'5.600"
$ws.Range('E33').Value = '  +1.59%  '
$ws.Range('D34').Value = "'3.570"
$ws.Range('E34').Value = '  +0.27%  '
$ws.Range('E35').Value = '  -0.38%  '
$ws.Range('D36').Value = "'0.06514"
$ws.Range('E36').Value = '  +6.44%  '
$ws.Range('D37').Value = "'0.02283"
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').Value = "'1.197"
$ws.Range('E38').Value = '  +1.41%  '
$ws.Range('D39').Value = "'0.5939"
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('E40').Value = '  -0.92%  '
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('D42').Value = "'7.871"
$ws.Range('E42').Value = '  -1.22%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = "'2.515"
$ws.Range('E43').Value = '  +5.83%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').Value = "'0.1842"
$ws.Range('E44').Value = '  -0.79%  '
$ws.Range('D45').Value = "'1.246"
$ws.Range('E45').Value = '  -2.44%  '
$ws.Range('D46').Value = "'12.41"
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = "'0.07493"
$ws.Range('E47').Value = '  -1.47%  '
$ws.Range('D48').Value = "'0.5546"
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('D49').Value = "'1.963"
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('D50').Value = "'118.30"
$ws.Range('E50').Value = '  +1.42%  '
$ws.Range('D51').Value = "'2.436"
$ws.Range('E51').Value = '  +0.99%  '
